$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D, shifting old D:K to F:M.
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy number formats/styles from the (now shifted) old "D" column (now F)
# into the two newly inserted blank columns so the new quarter columns
# look like the rest of the data (date format row, number format rows, etc).
# Only the rows that actually contained data before the edit are targeted,
# so we don't fabricate empty styled cells on label-only rows. Each
# contiguous block is pasted separately (a single multi-area paste would
# also stamp the in-between label rows).
$ws.Range("F7:F35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)

$ws.Range("F38:F77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)

$ws.Range("F80:F102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Populate the two new quarter columns (D = newest quarter, E = next) with
# the latest reported financial data.
$ws.Range("D7").Value2 = 43465
$ws.Range("E7").Value2 = 43373
$ws.Range("D8").Value2 = 7945000
$ws.Range("E8").Value2 = 8152000
$ws.Range("D9").Value2 = 4063000
$ws.Range("E9").Value2 = 4160000
$ws.Range("D10").Value2 = 3882000
$ws.Range("E10").Value2 = 3992000
$ws.Range("D12").Value2 = 437000
$ws.Range("E12").Value2 = 430000
$ws.Range("D13").Value2 = 0
$ws.Range("E13").Value2 = 0
$ws.Range("D14").Value2 = -17000
$ws.Range("E14").Value2 = 0
$ws.Range("D15").Value2 = 0
$ws.Range("E15").Value2 = 0
$ws.Range("D17").Value2 = 6162000
$ws.Range("E17").Value2 = 6136000
$ws.Range("D18").Value2 = 1783000
$ws.Range("E18").Value2 = 2016000
$ws.Range("D20").Value2 = 32000
$ws.Range("E20").Value2 = 34000
$ws.Range("D21").Value2 = 2186000
$ws.Range("E21").Value2 = 2405000
$ws.Range("D22").Value2 = 95000
$ws.Range("E22").Value2 = 85000
$ws.Range("D23").Value2 = 1720000
$ws.Range("E23").Value2 = 1965000
$ws.Range("D24").Value2 = 412000
$ws.Range("E24").Value2 = 419000
$ws.Range("D25").Value2 = 0
$ws.Range("E25").Value2 = 0
$ws.Range("D26").Value2 = 1308000
$ws.Range("E26").Value2 = 1546000
$ws.Range("D27").Value2 = 1306000
$ws.Range("E27").Value2 = 1543000
$ws.Range("D28").Value2 = 0
$ws.Range("E28").Value2 = 0
$ws.Range("D29").Value2 = 41000
$ws.Range("E29").Value2 = 0
$ws.Range("D30").Value2 = 0
$ws.Range("E30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("E31").Value2 = 0
$ws.Range("D32").Value2 = -32000
$ws.Range("E32").Value2 = -34000
$ws.Range("D33").Value2 = 1347000
$ws.Range("E33").Value2 = 1543000
$ws.Range("D34").Value2 = 0
$ws.Range("E34").Value2 = 0
$ws.Range("D35").Value2 = 1347000
$ws.Range("E35").Value2 = 1543000
$ws.Range("D38").Value2 = 43465
$ws.Range("E38").Value2 = 43373
$ws.Range("D41").Value2 = 2853000
$ws.Range("E41").Value2 = 3185000
$ws.Range("D42").Value2 = 380000
$ws.Range("E42").Value2 = 338000
$ws.Range("D43").Value2 = 5123000
$ws.Range("E43").Value2 = 5329000
$ws.Range("D44").Value2 = 4366000
$ws.Range("E44").Value2 = 4437000
$ws.Range("D45").Value2 = 987000
$ws.Range("E45").Value2 = 1130000
$ws.Range("D46").Value2 = 13709000
$ws.Range("E46").Value2 = 14419000
$ws.Range("D47:E47").Value2 = "NA"
$ws.Range("D48").Value2 = 8738000
$ws.Range("E48").Value2 = 8630000
$ws.Range("D49").Value2 = 12708000
$ws.Range("E49").Value2 = 12849000
$ws.Range("D50").Value2 = 0
$ws.Range("E50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("E51").Value2 = 0
$ws.Range("D52").Value2 = 1345000
$ws.Range("E52").Value2 = 1377000
$ws.Range("D53").Value2 = 0
$ws.Range("E53").Value2 = 0
$ws.Range("D54").Value2 = 36500000
$ws.Range("E54").Value2 = 37275000
$ws.Range("D57").Value2 = 2266000
$ws.Range("E57").Value2 = 2029000
$ws.Range("D58").Value2 = 1211000
$ws.Range("E58").Value2 = 1307000
$ws.Range("D59").Value2 = 3767000
$ws.Range("E59").Value2 = 4000000
$ws.Range("D60").Value2 = 7244000
$ws.Range("E60").Value2 = 7336000
$ws.Range("D61").Value2 = 13411000
$ws.Range("E61").Value2 = 13539000
$ws.Range("D62").Value2 = 5997000
$ws.Range("E62").Value2 = 6089000
$ws.Range("D63").Value2 = 0
$ws.Range("E63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("E64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("E65").Value2 = 0
$ws.Range("D66").Value2 = 26704000
$ws.Range("E66").Value2 = 27027000
$ws.Range("D68").Value2 = 0
$ws.Range("E68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("E69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("E70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("E71").Value2 = 0
$ws.Range("D72").Value2 = 40636000
$ws.Range("E72").Value2 = 40120000
$ws.Range("D73").Value2 = 0
$ws.Range("E73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("E74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("E75").Value2 = 0
$ws.Range("D76").Value2 = 9796000
$ws.Range("E76").Value2 = 10248000
$ws.Range("D77").Value2 = 0
$ws.Range("E77").Value2 = 0
$ws.Range("D80").Value2 = 43465
$ws.Range("E80").Value2 = 43373
$ws.Range("D81").Value2 = 1347000
$ws.Range("E81").Value2 = 1543000
$ws.Range("D83").Value2 = 371000
$ws.Range("E83").Value2 = 355000
$ws.Range("D84").Value2 = 0
$ws.Range("E84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("E85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("E86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("E87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("E88").Value2 = 0
$ws.Range("D89").Value2 = 2258000
$ws.Range("E89").Value2 = 2139000
$ws.Range("D91").Value2 = -531000
$ws.Range("E91").Value2 = -377000
$ws.Range("D92").Value2 = 0
$ws.Range("E92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("E93").Value2 = 0
$ws.Range("D94").Value2 = -416000
$ws.Range("E94").Value2 = -269000
$ws.Range("D96").Value2 = -787000
$ws.Range("E96").Value2 = -794000
$ws.Range("D97").Value2 = 0
$ws.Range("E97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("E98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("E99").Value2 = 0
$ws.Range("D100").Value2 = -2152000
$ws.Range("E100").Value2 = -1453000
$ws.Range("D101").Value2 = -22000
$ws.Range("E101").Value2 = -33000
$ws.Range("D102").Value2 = -332000
$ws.Range("E102").Value2 = 384000
